$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B for "Variável" (shifts old B "Ano" -> C, old C "Valor" -> D)
$ws.Columns("B").Insert()

# Insert a new column before E for "Colocação" (column E is currently blank after the first insert)
$ws.Columns("E").Insert()

# ----- Header row -----
$ws.Range("A1").Value = "Região"
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Ano"
$ws.Range("D1").Value = "Valor"
$ws.Range("E1").Value = "Colocação"

# ----- Column B: Variável (same text for every data row) -----
$variable = "Trabalho como origem na renda (%)"
$ws.Range("B2").Value = $variable
$ws.Range("B3").Value = $variable
$ws.Range("B4").Value = $variable
$ws.Range("B5").Value = $variable
$ws.Range("B6").Value = $variable
$ws.Range("B7").Value = $variable
$ws.Range("B8").Value = $variable
$ws.Range("B9").Value = $variable
$ws.Range("B10").Value = $variable

# ----- Column C: Ano, update date format from 31/12/19 to 31/12/2019 -----
$ws.Range("C2").Value = "31/12/2019"
$ws.Range("C3").Value = "31/12/2019"
$ws.Range("C4").Value = "31/12/2019"
$ws.Range("C5").Value = "31/12/2019"
$ws.Range("C6").Value = "31/12/2019"
$ws.Range("C7").Value = "31/12/2019"
$ws.Range("C8").Value = "31/12/2019"
$ws.Range("C9").Value = "31/12/2019"
$ws.Range("C10").Value = "31/12/2019"

# ----- Column D: Valor (values unchanged, carried over from old column C) -----
$ws.Range("D2").Value = 82.58002999999999
$ws.Range("D3").Value = 82.49911
$ws.Range("D4").Value = 80.19324
$ws.Range("D5").Value = 77.98804
$ws.Range("D6").Value = 77.80428999999999
$ws.Range("D7").Value = 76.63402000000001
$ws.Range("D8").Value = 66.29507
$ws.Range("D9").Value = 65.81395999999999
$ws.Range("D10").Value = 72.46505999999999

# ----- Column E: Colocação (ranking), only populated for rows 2-8 -----
$ws.Range("E2").Value = "1º"
$ws.Range("E3").Value = "2º"
$ws.Range("E4").Value = "3º"
$ws.Range("E5").Value = "4º"
$ws.Range("E6").Value = "5º"
$ws.Range("E7").Value = "6º"
$ws.Range("E8").Value = "26º"

# ----- Update the sheet dimension to reflect the new range -----
$ws.Range("A1:E10").Select()
